# Smoke test data cleanup + deployment changes
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # studiologin
$ws2 = $wb.Worksheets.Item(2)   # Organisation
$ws3 = $wb.Worksheets.Item(3)   # Process
$ws4 = $wb.Worksheets.Item(4)   # Entity
$ws5 = $wb.Worksheets.Item(5)   # dashboard

# ---------------------------------------------------------------------
# Sheet4 (Entity): rename the shared string "DrivingLeccDetail" -> "riskk"
# (B4 is the only cell referencing it directly by content)
# ---------------------------------------------------------------------
$ws4.Range("B4").Value2 = "riskk"

# ---------------------------------------------------------------------
# Sheet3 (Process): drop the Process101 block (rows 9-13), fix the
# Submodules cell, and update the selection / scroll position.
# ---------------------------------------------------------------------
$ws3.Range("C2").Value2 = "subMod1"
$ws3.Rows("9:13").Delete()

# ---------------------------------------------------------------------
# Sheet2 (Organisation): reshape the 3-column sparse layout into a dense
# 2-column Org/Head table.
# ---------------------------------------------------------------------
$ws2.Rows("2:7").Delete()
$ws2.Range("A1").Value2 = "Org"
$ws2.Range("B1").Value2 = "Head"
$ws2.Range("C1").Clear()
$ws2.Range("A2").Value2 = "CEO"
$ws2.Range("B2").Value2 = "CTO, CO, CFO"
$ws2.Range("A3").Value2 = "CTO"
$ws2.Range("B3").Value2 = "Engineering, Testing"
$ws2.Range("A4").Value2 = "CO"
$ws2.Range("B4").Value2 = "Admin, Functional"
$ws2.Range("A5").Value2 = "CFO"
$ws2.Range("B5").Value2 = "Finanacial,HR"
$ws2.Range("A6").Value2 = "save"
$ws2.Columns.Item(2).ColumnWidth = 18.5

# ---------------------------------------------------------------------
# Selections / active sheet. Touch the non-active sheets first so the
# last Activate()/Select() wins and leaves Organisation as the tab that
# is actually selected (tabSelected + workbook activeTab).
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D2").Select() | Out-Null

$ws3.Activate()
$ws3.Range("E4").Select() | Out-Null

$ws4.Activate()
$ws4.Range("B4").Select() | Out-Null

$ws5.Activate()
$ws5.Range("G19").Select() | Out-Null

$ws2.Activate()
$ws2.Range("B12").Select() | Out-Null
